# Update Excel data - 2024-11-22 04:46:16
# Refresh live crypto data across all three sheets: per-coin rows on
# "Top 50 Cryptocurrencies" (including a few rank swaps caused by the
# live re-sort by market cap), the "Top 5 by Market Cap" rollup, and the
# computed figures on "Summary".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws3 = $wb.Worksheets.Item("Summary")

# --- Sheet 1: Top 50 Cryptocurrencies ---
$ws1.Cells.Item(2, 1).Value = "Bitcoin"
$ws1.Cells.Item(2, 2).Value = "btc"
$ws1.Cells.Item(2, 3).Value = 99149
$ws1.Cells.Item(2, 4).Value = 1961769311574
$ws1.Cells.Item(2, 5).Value = 98540898109
$ws1.Cells.Item(2, 6).Value = 2.39324

$ws1.Cells.Item(3, 1).Value = "Ethereum"
$ws1.Cells.Item(3, 2).Value = "eth"
$ws1.Cells.Item(3, 3).Value = 3414.08
$ws1.Cells.Item(3, 4).Value = 411523313271
$ws1.Cells.Item(3, 5).Value = 55309814282
$ws1.Cells.Item(3, 6).Value = 8.86512

$ws1.Cells.Item(4, 1).Value = "Tether"
$ws1.Cells.Item(4, 2).Value = "usdt"
$ws1.Cells.Item(4, 3).Value = 1.002
$ws1.Cells.Item(4, 4).Value = 130924741778
$ws1.Cells.Item(4, 5).Value = 105189143813
$ws1.Cells.Item(4, 6).Value = -0.00751

$ws1.Cells.Item(5, 1).Value = "Solana"
$ws1.Cells.Item(5, 2).Value = "sol"
$ws1.Cells.Item(5, 3).Value = 261.82
$ws1.Cells.Item(5, 4).Value = 124291615827
$ws1.Cells.Item(5, 5).Value = 15188480023
$ws1.Cells.Item(5, 6).Value = 8.75882

$ws1.Cells.Item(6, 1).Value = "BNB"
$ws1.Cells.Item(6, 2).Value = "bnb"
$ws1.Cells.Item(6, 3).Value = 635.84
$ws1.Cells.Item(6, 4).Value = 92854694727
$ws1.Cells.Item(6, 5).Value = 2507906945
$ws1.Cells.Item(6, 6).Value = 4.3857

$ws1.Cells.Item(7, 1).Value = "XRP"
$ws1.Cells.Item(7, 2).Value = "xrp"
$ws1.Cells.Item(7, 3).Value = 1.4
$ws1.Cells.Item(7, 4).Value = 79880606217
$ws1.Cells.Item(7, 5).Value = 17813777637
$ws1.Cells.Item(7, 6).Value = 26.04755

$ws1.Cells.Item(8, 1).Value = "Dogecoin"
$ws1.Cells.Item(8, 2).Value = "doge"
$ws1.Cells.Item(8, 3).Value = 0.396717
$ws1.Cells.Item(8, 4).Value = 58283012455
$ws1.Cells.Item(8, 5).Value = 10244606719
$ws1.Cells.Item(8, 6).Value = 2.18157

$ws1.Cells.Item(9, 1).Value = "USDC"
$ws1.Cells.Item(9, 2).Value = "usdc"
$ws1.Cells.Item(9, 3).Value = 1
$ws1.Cells.Item(9, 4).Value = 38326204577
$ws1.Cells.Item(9, 5).Value = 15125096053
$ws1.Cells.Item(9, 6).Value = -0.09701

$ws1.Cells.Item(10, 1).Value = "Lido Staked Ether"
$ws1.Cells.Item(10, 2).Value = "steth"
$ws1.Cells.Item(10, 3).Value = 3412.1
$ws1.Cells.Item(10, 4).Value = 33430993756
$ws1.Cells.Item(10, 5).Value = 147143004
$ws1.Cells.Item(10, 6).Value = 8.99408

$ws1.Cells.Item(11, 1).Value = "Cardano"
$ws1.Cells.Item(11, 2).Value = "ada"
$ws1.Cells.Item(11, 3).Value = 0.892538
$ws1.Cells.Item(11, 4).Value = 31989342572
$ws1.Cells.Item(11, 5).Value = 2843145921
$ws1.Cells.Item(11, 6).Value = 12.53038

$ws1.Cells.Item(12, 1).Value = "TRON"
$ws1.Cells.Item(12, 2).Value = "trx"
$ws1.Cells.Item(12, 3).Value = 0.201062
$ws1.Cells.Item(12, 4).Value = 17366628805
$ws1.Cells.Item(12, 5).Value = 1099122463
$ws1.Cells.Item(12, 6).Value = 1.75848

$ws1.Cells.Item(13, 1).Value = "Avalanche"
$ws1.Cells.Item(13, 2).Value = "avax"
$ws1.Cells.Item(13, 3).Value = 36.57
$ws1.Cells.Item(13, 4).Value = 14963010321
$ws1.Cells.Item(13, 5).Value = 1052962960
$ws1.Cells.Item(13, 6).Value = 7.13683

$ws1.Cells.Item(14, 1).Value = "Shiba Inu"
$ws1.Cells.Item(14, 2).Value = "shib"
$ws1.Cells.Item(14, 3).Value = 0.00002514
$ws1.Cells.Item(14, 4).Value = 14817582167
$ws1.Cells.Item(14, 5).Value = 1622034870
$ws1.Cells.Item(14, 6).Value = 4.32907

$ws1.Cells.Item(15, 1).Value = "Wrapped stETH"
$ws1.Cells.Item(15, 2).Value = "wsteth"
$ws1.Cells.Item(15, 3).Value = 4010.15
$ws1.Cells.Item(15, 4).Value = 14476726400
$ws1.Cells.Item(15, 5).Value = 166791664
$ws1.Cells.Item(15, 6).Value = 8.59362

$ws1.Cells.Item(16, 1).Value = "Wrapped Bitcoin"
$ws1.Cells.Item(16, 2).Value = "wbtc"
$ws1.Cells.Item(16, 3).Value = 98631
$ws1.Cells.Item(16, 4).Value = 14410603666
$ws1.Cells.Item(16, 5).Value = 886179449
$ws1.Cells.Item(16, 6).Value = 2.50796

$ws1.Cells.Item(17, 1).Value = "Toncoin"
$ws1.Cells.Item(17, 2).Value = "ton"
$ws1.Cells.Item(17, 3).Value = 5.59
$ws1.Cells.Item(17, 4).Value = 14218996814
$ws1.Cells.Item(17, 5).Value = 630148511
$ws1.Cells.Item(17, 6).Value = 4.15114

$ws1.Cells.Item(18, 1).Value = "Sui"
$ws1.Cells.Item(18, 2).Value = "sui"
$ws1.Cells.Item(18, 3).Value = 3.64
$ws1.Cells.Item(18, 4).Value = 10374254368
$ws1.Cells.Item(18, 5).Value = 2412961752
$ws1.Cells.Item(18, 6).Value = 0.9307

$ws1.Cells.Item(19, 1).Value = "Bitcoin Cash"
$ws1.Cells.Item(19, 2).Value = "bch"
$ws1.Cells.Item(19, 3).Value = 496.39
$ws1.Cells.Item(19, 4).Value = 9814281450
$ws1.Cells.Item(19, 5).Value = 2241782063
$ws1.Cells.Item(19, 6).Value = 4.02434

$ws1.Cells.Item(20, 1).Value = "WETH"
$ws1.Cells.Item(20, 2).Value = "weth"
$ws1.Cells.Item(20, 3).Value = 3414.38
$ws1.Cells.Item(20, 4).Value = 9729290944
$ws1.Cells.Item(20, 5).Value = 1145242149
$ws1.Cells.Item(20, 6).Value = 9.15644

$ws1.Cells.Item(21, 1).Value = "Chainlink"
$ws1.Cells.Item(21, 2).Value = "link"
$ws1.Cells.Item(21, 3).Value = 15.37
$ws1.Cells.Item(21, 4).Value = 9636166601
$ws1.Cells.Item(21, 5).Value = 1249094715
$ws1.Cells.Item(21, 6).Value = 6.17806

$ws1.Cells.Item(22, 1).Value = "Pepe"
$ws1.Cells.Item(22, 2).Value = "pepe"
$ws1.Cells.Item(22, 3).Value = 0.00002149
$ws1.Cells.Item(22, 4).Value = 9079218772
$ws1.Cells.Item(22, 5).Value = 7071683603
$ws1.Cells.Item(22, 6).Value = 10.8927

$ws1.Cells.Item(23, 1).Value = "Polkadot"
$ws1.Cells.Item(23, 2).Value = "dot"
$ws1.Cells.Item(23, 3).Value = 6.26
$ws1.Cells.Item(23, 4).Value = 9008856710
$ws1.Cells.Item(23, 5).Value = 821046539
$ws1.Cells.Item(23, 6).Value = 9.63602

$ws1.Cells.Item(24, 1).Value = "Stellar"
$ws1.Cells.Item(24, 2).Value = "xlm"
$ws1.Cells.Item(24, 3).Value = 0.288277
$ws1.Cells.Item(24, 4).Value = 8632018977
$ws1.Cells.Item(24, 5).Value = 2350324462
$ws1.Cells.Item(24, 6).Value = 20.85285

$ws1.Cells.Item(25, 1).Value = "LEO Token"
$ws1.Cells.Item(25, 2).Value = "leo"
$ws1.Cells.Item(25, 3).Value = 8.76
$ws1.Cells.Item(25, 4).Value = 8097734319
$ws1.Cells.Item(25, 5).Value = 3467097
$ws1.Cells.Item(25, 6).Value = 2.69848

$ws1.Cells.Item(26, 1).Value = "NEAR Protocol"
$ws1.Cells.Item(26, 2).Value = "near"
$ws1.Cells.Item(26, 3).Value = 5.85
$ws1.Cells.Item(26, 4).Value = 7124913275
$ws1.Cells.Item(26, 5).Value = 1016659437
$ws1.Cells.Item(26, 6).Value = 4.39655

$ws1.Cells.Item(27, 1).Value = "Litecoin"
$ws1.Cells.Item(27, 2).Value = "ltc"
$ws1.Cells.Item(27, 3).Value = 91.07
$ws1.Cells.Item(27, 4).Value = 6847657317
$ws1.Cells.Item(27, 5).Value = 1465484425
$ws1.Cells.Item(27, 6).Value = 5.62763

$ws1.Cells.Item(28, 1).Value = "Aptos"
$ws1.Cells.Item(28, 2).Value = "apt"
$ws1.Cells.Item(28, 3).Value = 12.17
$ws1.Cells.Item(28, 4).Value = 6486163977
$ws1.Cells.Item(28, 5).Value = 892733968
$ws1.Cells.Item(28, 6).Value = 3.53089

$ws1.Cells.Item(29, 1).Value = "Wrapped eETH"
$ws1.Cells.Item(29, 2).Value = "weeth"
$ws1.Cells.Item(29, 3).Value = 3556.65
$ws1.Cells.Item(29, 4).Value = 6189609043
$ws1.Cells.Item(29, 5).Value = 100989166
$ws1.Cells.Item(29, 6).Value = 7.64083

$ws1.Cells.Item(30, 1).Value = "Uniswap"
$ws1.Cells.Item(30, 2).Value = "uni"
$ws1.Cells.Item(30, 3).Value = 9.5
$ws1.Cells.Item(30, 4).Value = 5704273647
$ws1.Cells.Item(30, 5).Value = 857998246
$ws1.Cells.Item(30, 6).Value = 8.0477

$ws1.Cells.Item(31, 1).Value = "Cronos"
$ws1.Cells.Item(31, 2).Value = "cro"
$ws1.Cells.Item(31, 3).Value = 0.193386
$ws1.Cells.Item(31, 4).Value = 5241773560
$ws1.Cells.Item(31, 5).Value = 116796609
$ws1.Cells.Item(31, 6).Value = 9.25628

$ws1.Cells.Item(32, 1).Value = "USDS"
$ws1.Cells.Item(32, 2).Value = "usds"
$ws1.Cells.Item(32, 3).Value = 1.001
$ws1.Cells.Item(32, 4).Value = 5240540218
$ws1.Cells.Item(32, 5).Value = 16261713
$ws1.Cells.Item(32, 6).Value = 0.19536

$ws1.Cells.Item(33, 1).Value = "Hedera"
$ws1.Cells.Item(33, 2).Value = "hbar"
$ws1.Cells.Item(33, 3).Value = 0.134935
$ws1.Cells.Item(33, 4).Value = 5109436515
$ws1.Cells.Item(33, 5).Value = 885897353
$ws1.Cells.Item(33, 6).Value = 6.79599

$ws1.Cells.Item(34, 1).Value = "Internet Computer"
$ws1.Cells.Item(34, 2).Value = "icp"
$ws1.Cells.Item(34, 3).Value = 9.77
$ws1.Cells.Item(34, 4).Value = 4634790671
$ws1.Cells.Item(34, 5).Value = 273253536
$ws1.Cells.Item(34, 6).Value = 7.39788

$ws1.Cells.Item(35, 1).Value = "Ethereum Classic"
$ws1.Cells.Item(35, 2).Value = "etc"
$ws1.Cells.Item(35, 3).Value = 28.22
$ws1.Cells.Item(35, 4).Value = 4223356412
$ws1.Cells.Item(35, 5).Value = 895949110
$ws1.Cells.Item(35, 6).Value = 7.09025

$ws1.Cells.Item(36, 1).Value = "Bonk"
$ws1.Cells.Item(36, 2).Value = "bonk"
$ws1.Cells.Item(36, 3).Value = 0.00005249
$ws1.Cells.Item(36, 4).Value = 3932666747
$ws1.Cells.Item(36, 5).Value = 1777055738
$ws1.Cells.Item(36, 6).Value = 4.15601

$ws1.Cells.Item(37, 1).Value = "Render"
$ws1.Cells.Item(37, 2).Value = "render"
$ws1.Cells.Item(37, 3).Value = 7.46
$ws1.Cells.Item(37, 4).Value = 3863154821
$ws1.Cells.Item(37, 5).Value = 444152804
$ws1.Cells.Item(37, 6).Value = 1.04288

$ws1.Cells.Item(38, 1).Value = "Kaspa"
$ws1.Cells.Item(38, 2).Value = "kas"
$ws1.Cells.Item(38, 3).Value = 0.150858
$ws1.Cells.Item(38, 4).Value = 3808781070
$ws1.Cells.Item(38, 5).Value = 154293525
$ws1.Cells.Item(38, 6).Value = -0.2823

$ws1.Cells.Item(39, 1).Value = "POL (ex-MATIC)"
$ws1.Cells.Item(39, 2).Value = "pol"
$ws1.Cells.Item(39, 3).Value = 0.477072
$ws1.Cells.Item(39, 4).Value = 3802839938
$ws1.Cells.Item(39, 5).Value = 486060217
$ws1.Cells.Item(39, 6).Value = 8.22479

$ws1.Cells.Item(40, 1).Value = "Bittensor"
$ws1.Cells.Item(40, 2).Value = "tao"
$ws1.Cells.Item(40, 3).Value = 510.62
$ws1.Cells.Item(40, 4).Value = 3768031715
$ws1.Cells.Item(40, 5).Value = 287161985
$ws1.Cells.Item(40, 6).Value = 4.16779

$ws1.Cells.Item(41, 1).Value = "Ethena USDe"
$ws1.Cells.Item(41, 2).Value = "usde"
$ws1.Cells.Item(41, 3).Value = 1.002
$ws1.Cells.Item(41, 4).Value = 3689043815
$ws1.Cells.Item(41, 5).Value = 228245435
$ws1.Cells.Item(41, 6).Value = -0.08268

$ws1.Cells.Item(42, 1).Value = "WhiteBIT Coin"
$ws1.Cells.Item(42, 2).Value = "wbt"
$ws1.Cells.Item(42, 3).Value = 24.84
$ws1.Cells.Item(42, 4).Value = 3580056600
$ws1.Cells.Item(42, 5).Value = 42172745
$ws1.Cells.Item(42, 6).Value = 2.81733

$ws1.Cells.Item(43, 1).Value = "MANTRA"
$ws1.Cells.Item(43, 2).Value = "om"
$ws1.Cells.Item(43, 3).Value = 3.84
$ws1.Cells.Item(43, 4).Value = 3462791802
$ws1.Cells.Item(43, 5).Value = 305531810
$ws1.Cells.Item(43, 6).Value = 5.79364

$ws1.Cells.Item(44, 1).Value = "Dai"
$ws1.Cells.Item(44, 2).Value = "dai"
$ws1.Cells.Item(44, 3).Value = 1.001
$ws1.Cells.Item(44, 4).Value = 3444417149
$ws1.Cells.Item(44, 5).Value = 160828033
$ws1.Cells.Item(44, 6).Value = -0.19239

$ws1.Cells.Item(45, 1).Value = "dogwifhat"
$ws1.Cells.Item(45, 2).Value = "wif"
$ws1.Cells.Item(45, 3).Value = 3.41
$ws1.Cells.Item(45, 4).Value = 3406645100
$ws1.Cells.Item(45, 5).Value = 1283305572
$ws1.Cells.Item(45, 6).Value = 5.93067

$ws1.Cells.Item(46, 1).Value = "Artificial Superintelligence Alliance"
$ws1.Cells.Item(46, 2).Value = "fet"
$ws1.Cells.Item(46, 3).Value = 1.29
$ws1.Cells.Item(46, 4).Value = 3367500054
$ws1.Cells.Item(46, 5).Value = 477414502
$ws1.Cells.Item(46, 6).Value = 4.08397

$ws1.Cells.Item(47, 1).Value = "Arbitrum"
$ws1.Cells.Item(47, 2).Value = "arb"
$ws1.Cells.Item(47, 3).Value = 0.794326
$ws1.Cells.Item(47, 4).Value = 3255030415
$ws1.Cells.Item(47, 5).Value = 1679348260
$ws1.Cells.Item(47, 6).Value = 12.9665

$ws1.Cells.Item(48, 1).Value = "Monero"
$ws1.Cells.Item(48, 2).Value = "xmr"
$ws1.Cells.Item(48, 3).Value = 160.74
$ws1.Cells.Item(48, 4).Value = 2965017866
$ws1.Cells.Item(48, 5).Value = 86598350
$ws1.Cells.Item(48, 6).Value = -0.70978

$ws1.Cells.Item(49, 1).Value = "Stacks"
$ws1.Cells.Item(49, 2).Value = "stx"
$ws1.Cells.Item(49, 3).Value = 1.96
$ws1.Cells.Item(49, 4).Value = 2941212034
$ws1.Cells.Item(49, 5).Value = 390052677
$ws1.Cells.Item(49, 6).Value = 2.29854

$ws1.Cells.Item(50, 1).Value = "Filecoin"
$ws1.Cells.Item(50, 2).Value = "fil"
$ws1.Cells.Item(50, 3).Value = 4.73
$ws1.Cells.Item(50, 4).Value = 2842098539
$ws1.Cells.Item(50, 5).Value = 594063403
$ws1.Cells.Item(50, 6).Value = 8.51863

$ws1.Cells.Item(51, 1).Value = "OKB"
$ws1.Cells.Item(51, 2).Value = "okb"
$ws1.Cells.Item(51, 3).Value = 46.73
$ws1.Cells.Item(51, 4).Value = 2804621854
$ws1.Cells.Item(51, 5).Value = 20202917
$ws1.Cells.Item(51, 6).Value = 5.99002

# --- Sheet 2: Top 5 by Market Cap ---
$ws2.Cells.Item(2, 2).Value = 1961769311574
$ws2.Cells.Item(3, 2).Value = 411523313271
$ws2.Cells.Item(4, 2).Value = 130924741778
$ws2.Cells.Item(5, 2).Value = 124291615827
$ws2.Cells.Item(6, 2).Value = 92854694727

# --- Sheet 3: Summary ---
# The average-price cell holds a literal "$..." text label (not a real
# currency number) in the source file. Assigning a "$"-led string straight
# to .Value lets Excel's smart-text heuristics reinterpret it as a Currency
# number, so force the cell to Text first, then restore the default style
# (the source cell carries no explicit style) once the text is in place.
$ws3.Cells.Item(2, 2).NumberFormat = "@"
$ws3.Cells.Item(2, 2).Value = "$4359.79"
$ws3.Cells.Item(2, 2).Style = "Normal"
$ws3.Cells.Item(3, 2).Value = "XRP (26.05%)"
$ws3.Cells.Item(4, 2).Value = "Monero (-0.71%)"

